$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the header style used by
# the other header cells (B1:E1) -- copy formatting from E1 so it reuses
# the existing style instead of creating a new one.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the time_taken values for each data row (unstyled, like the
# other data cells).
$ws.Range("F2").Value = "2021-10-05 13:42:10.268447"
$ws.Range("F3").Value = "2021-10-05 13:42:10.268459"
$ws.Range("F4").Value = "2021-10-05 13:42:10.268463"
$ws.Range("F5").Value = "2021-10-05 13:42:10.268466"
$ws.Range("F6").Value = "2021-10-05 13:42:10.268470"
$ws.Range("F7").Value = "2021-10-05 13:42:10.268473"
$ws.Range("F8").Value = "2021-10-05 13:42:10.268476"
$ws.Range("F9").Value = "2021-10-05 13:42:10.268479"
